$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The order-number and amount columns hold numeric-looking text values, so
# force just those cells to Text format before writing them, keeping them
# as shared strings instead of being auto-converted to real numbers.
$ws.Range("A8:A10").NumberFormat = "@"
$ws.Range("D8:D10").NumberFormat = "@"

# Update period text (merged cell A2:G2)
$ws.Range("A2").Value = "Период: 2023-11-01 - 2023-11-30"

# Update data row 8 (first data row)
$ws.Cells.Item(8, 1).Value = "70004"
$ws.Cells.Item(8, 2).Value = "Сбербанк"
$ws.Cells.Item(8, 3).Value = "лебенков"
$ws.Cells.Item(8, 4).Value = "11300"
$ws.Cells.Item(8, 5).Value = "01.11.2023 00:00:00"
$ws.Cells.Item(8, 6).Value = "Поступление"
$ws.Cells.Item(8, 7).Value = "Докторская, Волковыское"

# Update data row 9
$ws.Cells.Item(9, 1).Value = "70005"
$ws.Cells.Item(9, 2).Value = "Сбербанк"
$ws.Cells.Item(9, 3).Value = "лебенков"
$ws.Cells.Item(9, 4).Value = "5000"
$ws.Cells.Item(9, 5).Value = "01.11.2023 00:00:00"
$ws.Cells.Item(9, 6).Value = "Выбытие"
$ws.Cells.Item(9, 7).Value = "Докторская"

# Update data row 10
$ws.Cells.Item(10, 1).Value = "70006"
$ws.Cells.Item(10, 2).Value = "Сбербанк"
$ws.Cells.Item(10, 3).Value = "лебенков"
$ws.Cells.Item(10, 4).Value = "5000"
$ws.Cells.Item(10, 5).Value = "01.11.2023 00:00:00"
$ws.Cells.Item(10, 6).Value = "Выбытие"
$ws.Cells.Item(10, 7).Value = "Докторская"

# Remove the now-obsolete row 11 entirely
$ws.Rows.Item(11).Delete()

# Widen column G to fit the new text
$ws.Columns.Item(7).ColumnWidth = 24.33
